$d = $word.ActiveDocument

# The template has a "Success story from this quarter" heading that
# introduces the `{{ success_story }}` merge field (inside the
# "{%p if success_story %} ... {%p endif %}" block). Rename it to
# "Biggest win this quarter".
$d.Content.Find.Execute("Success story from this quarter", $true, $true, $false, $false, $false,
                         $true, 1, $false, "Biggest win this quarter", 2) | Out-Null
